$wb = $excel.ActiveWorkbook

# Fix the typo in the "controller" sheet's param name: paralel_controlers -> parallel_controllers
$controllerSheet = $wb.Worksheets.Item("controller")
$controllerSheet.Range("A4").Value = "parallel_controllers"

# Change active sheet to "controller" (activeTab="1") and update selections.

# battery: no longer the selected tab, but selection stays at A13
$batterySheet = $wb.Worksheets.Item("battery")
$batterySheet.Range("A13").Select()

# controller: becomes the selected tab, selection moves to A5
$controllerSheet.Range("A5").Select()
$controllerSheet.Activate()

# panels: selection moves to A7 (tab not changed)
$panelsSheet = $wb.Worksheets.Item("panels")
$panelsSheet.Range("A7").Select()

# Re-activate controller sheet last so it ends up as the active/selected tab
$controllerSheet.Activate()
